$d = $word.ActiveDocument

$ids = @("p019v_1", "p019v_2", "p019v_3")
foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $null = $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
}
